$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 83 (shifts the existing rows 83:90 down to 84:91),
# keeping the alphabetically-sorted list of islands in sync ("The Reaper's
# Hideout" sorts between "The North Star Seapost" and "The Spoils of
# Plenty Store").
$ws.Rows("83:83").Insert()

$ws.Range("B83").Value = "I-12"
$ws.Range("A83").Value = "The Reaper’s Hideout"
$ws.Range("C83").Value = "The Reaper’s Hideout"

# Re-apply the sort over the now one-row-bigger table so the sheet's
# sortState/sortCondition references extend to the new last row (A2:C91).
$sortRange = $ws.Range("A2:C91")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A91"))
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# Keep the view/selection in sync with the saved state.
$ws.Range("A83").Select()
$excel.ActiveWindow.ScrollRow = 74
